$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the registration data row (row 2) with new test data
$ws.Range("A2").Value = "Sandeep"
$ws.Range("B2").Value = "sandeep@gmail.com"
$ws.Range("C2").Value = 7817008251
$ws.Range("D2").Value = "Vadodara"
$ws.Range("F2").Value = "Sandeep@123"
$ws.Range("G2").Value = "Sandeep@123"

# Move the active selection to B2
$ws.Range("B2").Select()
